$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range('D2').Value = '27.741.35'
$ws.Range('E2').Value = '  +1.32%  '

$ws.Range('D3').Value = '1.646.19'
$ws.Range('E3').Value = '  -0.55%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = '213.55'
$ws.Range('E5').Value = '  +0.11%  '

$ws.Range('E6').Value = '  +3.64%  '

$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.01%  '

$ws.Range('D8').Value = '23.13'
$ws.Range('E8').Value = '  -1.51%  '

$ws.Range('E9').Value = '  +0.13%  '

$ws.Range('E10').Value = '  +0.05%  '

$ws.Range('D11').Value = '0.0892'
$ws.Range('E11').Value = '  +2.00%  '

$ws.Range('D12').Value = '1.878.78'
$ws.Range('E12').Value = '  -0.57%  '

$ws.Range('D13').Value = '1.660.45'
$ws.Range('E13').Value = '  +0.27%  '

$ws.Range('E14').Value = '  -0.85%  '

$ws.Range('E15').Value = '  -1.21%  '

$ws.Range('D16').Value = '64.34'
$ws.Range('E16').Value = '  -1.80%  '

$ws.Range('D17').Value = '27.722.04'

$ws.Range('D18').Value = '232.11'
$ws.Range('E18').Value = '  +0.16%  '

$ws.Range('E19').Value = '  +0.15%  '

$ws.Range('E20').Value = '  +3.29%  '

$ws.Range('E21').Value = '  +0.06%  '

$ws.Range('E22').Value = '  -0.76%  '

$ws.Range('D23').Value = '10.09'
$ws.Range('E23').Value = '  +7.34%  '

$ws.Range('E24').Value = '  -2.94%  '

$ws.Range('D25').Value = '149.63'
$ws.Range('E25').Value = '  +1.34%  '

$ws.Range('D26').Value = '6.98'
$ws.Range('E26').Value = '  -1.61%  '

$ws.Range('E27').Value = '  +1.09%  '

$ws.Range('D28').Value = '15.69'
$ws.Range('E28').Value = '  -1.09%  '

$ws.Range('E29').Value = '  -0.04%  '

$ws.Range('E30').Value = '  +0.09%  '

$ws.Range('E31').Value = '  -2.07%  '

$ws.Range('D32').Value = '3.30'
$ws.Range('E32').Value = '  +0.37%  '

$ws.Range('E33').Value = '  +1.70%  '

$ws.Range('D34').Value = '1.445.11'
$ws.Range('E34').Value = '  +1.83%  '

$ws.Range('D35').Value = '1.59'
$ws.Range('E35').Value = '  +2.02%  '

$ws.Range('E36').Value = '  -0.98%  '

$ws.Range('E37').Value = '  +0.29%  '

$ws.Range('D38').Value = '0.885'
$ws.Range('E38').Value = '  -2.46%  '

$ws.Range('E39').Value = '  -0.49%  '

$ws.Range('D40').Value = '0.895'
$ws.Range('E40').Value = '  +11.86%  '

$ws.Range('E41').Value = '  -0.94%  '

$ws.Range('E42').Value = '  +0.03%  '

$ws.Range('E43').Value = '  +2.65%  '

$ws.Range('E44').Value = '  -0.58%  '

$ws.Range('E45').Value = '  +1.74%  '

$ws.Range('D46').Value = '65.92'
$ws.Range('E46').Value = '  +1.55%  '

$ws.Range('E47').Value = '  -0.54%  '

$ws.Range('E48').Value = '  +1.80%  '

$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '85.84'
$ws.Range('E49').Value = '  -2.36%  '

$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0106'
$ws.Range('E50').Value = '  -0.59%  '

$ws.Range('E51').Value = '  -1.70%  '
